# Insert a new data row at row 824 (pushes existing rows 824-904 down to 825-905)
# and populate it with the new "Ajo" (garlic) price record for Vega Modelo de Temuco.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(824).Insert()

$ws.Cells.Item(824, 1).Value  = 10
$ws.Cells.Item(824, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(824, 3).Value  = "La Araucanía"
$ws.Cells.Item(824, 4).Value  = 44946
$ws.Cells.Item(824, 5).Value  = 9
$ws.Cells.Item(824, 6).Value  = 100112003
$ws.Cells.Item(824, 7).Value  = "Ajo"
$ws.Cells.Item(824, 8).Value  = "Chino"
$ws.Cells.Item(824, 9).Value  = "Primera"
$ws.Cells.Item(824, 10).Value = 125
$ws.Cells.Item(824, 11).Value = 18000
$ws.Cells.Item(824, 12).Value = 18000
$ws.Cells.Item(824, 13).Value = 18000
$ws.Cells.Item(824, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(824, 15).Value = "China"
$ws.Cells.Item(824, 16).Value = 1800
$ws.Cells.Item(824, 17).Value = 10
$ws.Cells.Item(824, 18).Value = "Hortaliza"
